$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1608506101308365
$ws.Range("D2").Value = 0.165341538105249
$ws.Range("E2").Value = 0.1310240930147621
$ws.Range("F2").Value = 1.39207083062135
$ws.Range("G2").Value = 0.002423139866598334
$ws.Range("I2").Value = 0.3688388074871369
$ws.Range("J2").Value = 0.137291349773351
$ws.Range("M2").Value = 0.7647322773308929
$ws.Range("N2").Value = 1.315290495938342
$ws.Range("O2").Value = 3.500899659374511
$ws.Range("B3").Value = 0.1423618748002156
$ws.Range("D3").Value = 0.1664409464868015
$ws.Range("E3").Value = 0.1325557332136285
$ws.Range("F3").Value = 1.369494441465534
$ws.Range("G3").Value = 0.002427163182120056
$ws.Range("I3").Value = 0.3694989586157469
$ws.Range("J3").Value = 0.1394176267524401
$ws.Range("M3").Value = 0.694572930943977
$ws.Range("N3").Value = 1.278227656929744
$ws.Range("O3").Value = 3.419403185468127
$ws.Range("B4").Value = 0.1309860365615094
$ws.Range("D4").Value = 0.1671633826388614
$ws.Range("E4").Value = 0.1335555536175983
$ws.Range("F4").Value = 1.356530393248875
$ws.Range("G4").Value = 0.002429766796046149
$ws.Range("I4").Value = 0.3701280429628753
$ws.Range("J4").Value = 0.1408039308630835
$ws.Range("M4").Value = 0.6515058652062464
$ws.Range("N4").Value = 1.25589235021036
$ws.Range("O4").Value = 3.371694907330408
$ws.Range("B5").Value = 0.126344764268211
$ws.Range("D5").Value = 0.1674697133404095
$ws.Range("E5").Value = 0.1339779170882993
$ws.Range("F5").Value = 1.351472740341421
$ws.Range("G5").Value = 0.002430861412835248
$ws.Range("I5").Value = 0.3704408556211405
$ws.Range("J5").Value = 0.1413891387838753
$ws.Range("M5").Value = 0.6339595954394071
$ws.Range("N5").Value = 1.246897933924785
$ws.Range("O5").Value = 3.352837886300989
$ws.Range("B6").Value = 0.1255737634324703
$ws.Range("D6").Value = 0.1675213004099145
$ws.Range("E6").Value = 0.13404895172817
$ws.Range("F6").Value = 1.350646514302355
$ws.Range("G6").Value = 0.002431045207092771
$ws.Range("I6").Value = 0.3704962133840652
$ws.Range("J6").Value = 0.1414875359520362
$ws.Range("M6").Value = 0.6310463201532457
$ws.Range("N6").Value = 1.245410947098776
$ws.Range("O6").Value = 3.349741946757064
$ws.Range("B7").Value = 0.130923464544864
$ws.Range("D7").Value = 0.1671674655894435
$ws.Range("E7").Value = 0.1335611893204973
$ws.Range("F7").Value = 1.356461272392721
$ws.Range("G7").Value = 0.002429781422123536
$ws.Range("I7").Value = 0.3701320328202478
$ws.Range("J7").Value = 0.1408117411181395
$ws.Range("M7").Value = 0.6512692127499236
$ws.Range("N7").Value = 1.255770611591686
$ws.Range("O7").Value = 3.371438230546886
$ws.Range("B8").Value = 0.1544808599800547
$ws.Range("D8").Value = 0.1657107858862918
$ws.Range("E8").Value = 0.1315398739671991
$ws.Range("F8").Value = 1.384099813618988
$ws.Range("G8").Value = 0.002424499508277434
$ws.Range("I8").Value = 0.3690201023240256
$ws.Range("J8").Value = 0.1380077099591066
$ws.Range("M8").Value = 0.740539693841157
$ws.Range("N8").Value = 1.302424603267468
$ws.Range("O8").Value = 3.472314880006991
$ws.Range("B9").Value = 0.2004723384094405
$ws.Range("D9").Value = 0.1632296546840575
$ws.Range("E9").Value = 0.1280474675409247
$ws.Range("F9").Value = 1.445450909538621
$ws.Range("G9").Value = 0.002415194291100224
$ws.Range("I9").Value = 0.3686070059194329
$ws.Range("J9").Value = 0.1331511114278943
$ws.Range("M9").Value = 0.9156457298577863
$ws.Range("N9").Value = 1.397199755681953
$ws.Range("O9").Value = 3.688717252645858
$ws.Range("B10").Value = 0.2341185443242466
$ws.Range("D10").Value = 0.1616347364270538
$ws.Range("E10").Value = 0.1257691517323929
$ws.Range("F10").Value = 1.494930378922788
$ws.Range("G10").Value = 0.002408992541847799
$ws.Range("I10").Value = 0.3693709900524276
$ws.Range("J10").Value = 0.129976084875878
$ws.Range("M10").Value = 1.0442827736547
$ws.Range("N10").Value = 1.468766788754976
$ws.Range("O10").Value = 3.859182702852991
$ws.Range("B11").Value = 0.2493900399978486
$ws.Range("D11").Value = 0.1609584805111623
$ws.Range("E11").Value = 0.1247951675443533
$ws.Range("F11").Value = 1.518406494733796
$ws.Range("G11").Value = 0.002406307575935956
$ws.Range("I11").Value = 0.3699482995282395
$ws.Range("J11").Value = 0.1286174049579625
$ws.Range("M11").Value = 1.102792171155912
$ws.Range("N11").Value = 1.501731095903466
$ws.Range("O11").Value = 3.939256400993941
$ws.Range("B12").Value = 0.2551676310247331
$ws.Range("D12").Value = 0.1607094727261966
$ws.Range("E12").Value = 0.1244353290169276
$ws.Range("F12").Value = 1.527436125295097
$ws.Range("G12").Value = 0.002405310328937506
$ws.Range("I12").Value = 0.3701997739572604
$ws.Range("J12").Value = 0.1281152591509791
$ws.Range("M12").Value = 1.124946019625625
$ws.Range("N12").Value = 1.514271206361627
$ws.Range("O12").Value = 3.969944000119824
$ws.Range("B13").Value = 0.2539235705150702
$ws.Range("D13").Value = 0.160762786486579
$ws.Range("E13").Value = 0.1245124267657864
$ws.Range("F13").Value = 1.525485208627416
$ws.Range("G13").Value = 0.002405524238396722
$ws.Range("I13").Value = 0.3701441558488057
$ws.Range("J13").Value = 0.1282228550302378
$ws.Range("M13").Value = 1.120174917416236
$ws.Range("N13").Value = 1.511567946144766
$ws.Range("O13").Value = 3.963318593826557
$ws.Range("B14").Value = 0.2498654765765309
$ws.Range("D14").Value = 0.1609378527446772
$ws.Range("E14").Value = 0.1247653832165803
$ws.Range("F14").Value = 1.519146563993672
$ws.Range("G14").Value = 0.002406225142005116
$ws.Range("I14").Value = 0.3699683309495185
$ws.Range("J14").Value = 0.1285758452440913
$ws.Range("M14").Value = 1.104614836394006
$ws.Range("N14").Value = 1.502761639323694
$ws.Range("O14").Value = 3.941773753654445
$ws.Range("B15").Value = 0.2473790610428068
$ws.Range("D15").Value = 0.1610460070244848
$ws.Range("E15").Value = 0.1249214971778185
$ws.Range("F15").Value = 1.515282176276287
$ws.Range("G15").Value = 0.002406657000715748
$ws.Range("I15").Value = 0.369864907110582
$ws.Range("J15").Value = 0.1287936723416028
$ws.Range("M15").Value = 1.095083490496421
$ws.Range("N15").Value = 1.497374934230095
$ws.Range("O15").Value = 3.928624562740481
$ws.Range("B16").Value = 0.2331197862868066
$ws.Range("D16").Value = 0.1616799220619267
$ws.Range("E16").Value = 0.1258340610843653
$ws.Range("F16").Value = 1.493415681330731
$ws.Range("G16").Value = 0.002409170743682937
$ws.Range("I16").Value = 0.369337869583827
$ws.Range("J16").Value = 0.1300666047945906
$ws.Range("M16").Value = 1.040458777898792
$ws.Range("N16").Value = 1.466620582293586
$ws.Range("O16").Value = 3.854000731416534
$ws.Range("B17").Value = 0.2243630646781298
$ws.Range("D17").Value = 0.1620814225620499
$ws.Range("E17").Value = 0.12640988741264
$ws.Range("F17").Value = 1.480249550992426
$ws.Range("G17").Value = 0.002410747666532238
$ws.Range("I17").Value = 0.3690732694821399
$ws.Range("J17").Value = 0.1308694741860954
$ws.Range("M17").Value = 1.006945315608633
$ws.Range("N17").Value = 1.447857278716867
$ws.Range("O17").Value = 3.808870275562413
$ws.Range("B18").Value = 0.2193232256385045
$ws.Range("D18").Value = 0.1623169941711033
$ws.Range("E18").Value = 0.126746963360441
$ws.Range("F18").Value = 1.472767786077213
$ws.Range("G18").Value = 0.002411667500512447
$ws.Range("I18").Value = 0.3689427161766332
$ws.Range("J18").Value = 0.1313393233522362
$ws.Range("M18").Value = 0.9876685546359596
$ws.Range("N18").Value = 1.437103620934352
$ws.Range("O18").Value = 3.78315027476367
$ws.Range("B19").Value = 0.2176162863111415
$ws.Range("D19").Value = 0.1623975517798453
$ws.Range("E19").Value = 0.1268621003580117
$ws.Range("F19").Value = 1.470250208298438
$ws.Range("G19").Value = 0.002411981147202023
$ws.Range("I19").Value = 0.3689022348919373
$ws.Range("J19").Value = 0.1314997893032785
$ws.Range("M19").Value = 0.9811416905627794
$ws.Range("N19").Value = 1.433469273431228
$ws.Range("O19").Value = 3.774482717041337
$ws.Range("B20").Value = 0.2252955675920703
$ws.Range("D20").Value = 0.1620382020928002
$ws.Range("E20").Value = 0.1263479815127386
$ws.Range("F20").Value = 1.481641681016953
$ws.Range("G20").Value = 0.002410578473392559
$ws.Range("I20").Value = 0.3690991985135739
$ws.Range("J20").Value = 0.130783172912361
$ws.Range("M20").Value = 1.01051296248491
$ws.Range("N20").Value = 1.449850690492553
$ws.Range("O20").Value = 3.813649860142959
$ws.Range("B21").Value = 0.2510575871238245
$ws.Range("D21").Value = 0.1608862395958504
$ws.Range("E21").Value = 0.124690839756366
$ws.Range("F21").Value = 1.52100458049722
$ws.Range("G21").Value = 0.002406018741398455
$ws.Range("I21").Value = 0.3700190846094245
$ws.Range("J21").Value = 0.1284718277849137
$ws.Range("M21").Value = 1.109185283824587
$ws.Range("N21").Value = 1.505346722983461
$ws.Range("O21").Value = 3.948092062732769
$ws.Range("B22").Value = 0.2678629116925322
$ws.Range("D22").Value = 0.1601746029243643
$ws.Range("E22").Value = 0.1236601965146127
$ws.Range("F22").Value = 1.547545269998778
$ws.Range("G22").Value = 0.002403152255346678
$ws.Range("I22").Value = 0.3708117330353957
$ws.Range("J22").Value = 0.1270332735264752
$ws.Range("M22").Value = 1.173658861539224
$ws.Range("N22").Value = 1.541949582647874
$ws.Range("O22").Value = 4.038088693940892
$ws.Range("B23").Value = 0.2588966393191185
$ws.Range("D23").Value = 0.1605506471147429
$ws.Range("E23").Value = 0.124205472639189
$ws.Range("F23").Value = 1.533305258489463
$ws.Range("G23").Value = 0.0024046717949473
$ws.Range("I23").Value = 0.3703712232799674
$ws.Range("J23").Value = 0.1277944518767242
$ws.Range("M23").Value = 1.1392498219342
$ws.Range("N23").Value = 1.522383951809132
$ws.Range("O23").Value = 3.989860193946072
$ws.Range("B24").Value = 0.2248740000687235
$ws.Range("D24").Value = 0.1620577272908292
$ws.Range("E24").Value = 0.1263759504017381
$ws.Range("F24").Value = 1.481012026153238
$ws.Range("G24").Value = 0.002410654924317181
$ws.Range("I24").Value = 0.3690874088017395
$ws.Range("J24").Value = 0.1308221639648171
$ws.Range("M24").Value = 1.008900058597149
$ws.Range("N24").Value = 1.448949364212353
$ws.Range("O24").Value = 3.811488305817818
$ws.Range("B25").Value = 0.1880542968757766
$ws.Range("D25").Value = 0.1638607687199176
$ws.Range("E25").Value = 0.1289417722191022
$ws.Range("F25").Value = 1.428083535387444
$ws.Range("G25").Value = 0.002417599626549556
$ws.Range("I25").Value = 0.368530584173584
$ws.Range("J25").Value = 0.1343960367203003
$ws.Range("M25").Value = 0.868274248927051
$ws.Range("N25").Value = 1.371215769127332
$ws.Range("O25").Value = 3.628170192807715
